$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.527.03'
$ws.Range("E2").Value = '  -0.31%  '

$ws.Range("D3").Value = '1.714.01'
$ws.Range("E3").Value = '  -1.44%  '

$ws.Range("E4").Value = '  +0.10%  '

$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = '240.68'
$r.Style = "Normal"
$ws.Range("E5").Value = '  -2.29%  '

$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = '1.001'
$r.Style = "Normal"
$ws.Range("E6").Value = '  +0.08%  '

$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = '0.4907'
$r.Style = "Normal"
$ws.Range("E7").Value = '  -0.91%  '

$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = '0.2592'
$r.Style = "Normal"
$ws.Range("E8").Value = '  -3.15%  '

$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = '0.06191'
$r.Style = "Normal"
$ws.Range("E9").Value = '  -1.27%  '

$ws.Range("D10").Value = '1.732.81'
$ws.Range("E10").Value = '  -0.40%  '

$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = '0.06993'
$r.Style = "Normal"
$ws.Range("E11").Value = '  -0.73%  '

$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = '15.67'
$r.Style = "Normal"
$ws.Range("E12").Value = '  -0.50%  '

$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = '0.6046'
$r.Style = "Normal"
$ws.Range("E13").Value = '  -1.52%  '

$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = '4.459'
$r.Style = "Normal"
$ws.Range("E14").Value = '  -2.83%  '

$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = '76.61'
$r.Style = "Normal"
$ws.Range("E15").Value = '  -1.80%  '

$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = '1.000'
$r.Style = "Normal"
$ws.Range("E16").Value = '  +0.05%  '

$ws.Range("D17").Value = '26.422.47'
$ws.Range("E17").Value = '  -0.74%  '

$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = '1.001'
$r.Style = "Normal"
$ws.Range("E18").Value = '  +0.15%  '

$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = '0.000007132'
$r.Style = "Normal"
$ws.Range("E19").Value = '  -1.84%  '

$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = '11.32'
$r.Style = "Normal"
$ws.Range("E20").Value = '  -2.11%  '

$ws.Range("D21").Value = '1.952.25'
$ws.Range("E21").Value = '  -1.00%  '

$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = '4.404'
$r.Style = "Normal"
$ws.Range("E22").Value = '  -3.54%  '

$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = '8.478'
$r.Style = "Normal"
$ws.Range("E23").Value = '  -2.83%  '

$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = '5.073'
$r.Style = "Normal"
$ws.Range("E24").Value = '  -3.80%  '

$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = '137.76'
$r.Style = "Normal"
$ws.Range("E25").Value = '  -0.78%  '

$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = '15.24'
$r.Style = "Normal"
$ws.Range("E26").Value = '  -1.06%  '

$ws.Range("E27").Value = '  +1.11%  '

$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = '1.739'
$r.Style = "Normal"
$ws.Range("E28").Value = '  -0.94%  '

$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = '105.86'
$r.Style = "Normal"

$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = '3.900'
$r.Style = "Normal"
$ws.Range("E30").Value = '  -3.23%  '

$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = '0.07928'
$r.Style = "Normal"
$ws.Range("E31").Value = '  -1.44%  '

$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = '3.626'
$r.Style = "Normal"
$ws.Range("E32").Value = '  -2.86%  '

$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = '0.04498'
$r.Style = "Normal"
$ws.Range("E33").Value = '  -2.69%  '

$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = '2.641'
$r.Style = "Normal"
$ws.Range("E34").Value = '  +1.15%  '

$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = '0.9967'
$r.Style = "Normal"
$ws.Range("E35").Value = '  -1.75%  '

$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = '0.6236'
$r.Style = "Normal"
$ws.Range("E36").Value = '  -2.51%  '

$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = '0.9406'
$r.Style = "Normal"
$ws.Range("E37").Value = '  +3.87%  '

$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = '1.998'
$r.Style = "Normal"
$ws.Range("E38").Value = '  -3.44%  '

$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = '2.408'
$r.Style = "Normal"
$ws.Range("E39").Value = '  -0.74%  '

$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = '1.001'
$r.Style = "Normal"
$ws.Range("E40").Value = '  -0.16%  '

$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = '0.01498'
$r.Style = "Normal"
$ws.Range("E41").Value = '  -0.48%  '

$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = '99.12'
$r.Style = "Normal"
$ws.Range("E42").Value = '  -2.68%  '

$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = '5.504'
$r.Style = "Normal"
$ws.Range("E43").Value = '  +1.14%  '

$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = '0.3828'
$r.Style = "Normal"
$ws.Range("E44").Value = '  -2.68%  '

$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = '6.914'
$r.Style = "Normal"
$ws.Range("E45").Value = '  +0.80%  '

$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = '0.1151'
$r.Style = "Normal"
$ws.Range("E46").Value = '  -2.65%  '

$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = '0.05372'
$r.Style = "Normal"
$ws.Range("E47").Value = '  -0.34%  '

$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = '7.761'
$r.Style = "Normal"
$ws.Range("E48").Value = '  -0.10%  '

$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = '30.22'
$r.Style = "Normal"
$ws.Range("E49").Value = '  -1.44%  '

$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = '51.36'
$r.Style = "Normal"
$ws.Range("E50").Value = '  -0.98%  '

$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = '1.217'
$r.Style = "Normal"
$ws.Range("E51").Value = '  -2.99%  '
